# Weekly data refresh: two new price-report rows are published for
# "Terminal Hortofrutícola Agro Chillán - Alcachofa" and inserted at the
# top of the data block (right after the header row), pushing the
# existing records down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 7 (the first data row), shifting the
# previously-existing rows 7..91 down to become rows 9..93.
$ws.Rows("7:8").Insert()

# New row 7
$ws.Range("A7").Value = 7
$ws.Range("B7").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C7").Value = "Ñuble"
$ws.Range("D7").Value = 45092
$ws.Range("E7").Value = 16
$ws.Range("F7").Value = 100112013
$ws.Range("G7").Value = "Alcachofa"
$ws.Range("H7").Value = "Argentina(o)"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 40
$ws.Range("K7").Value = 16000
$ws.Range("L7").Value = 17000
$ws.Range("M7").Value = 16500
$ws.Range("N7").Value = "$/caja 50 unidades"
$ws.Range("O7").Value = "Provincia de Limarí"
$ws.Range("P7").Value = 330
$ws.Range("Q7").Value = 50
$ws.Range("R7").Value = "Hortaliza"

# New row 8
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C8").Value = "Ñuble"
$ws.Range("D8").Value = 45092
$ws.Range("E8").Value = 16
$ws.Range("F8").Value = 100112013
$ws.Range("G8").Value = "Alcachofa"
$ws.Range("H8").Value = "Española"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 60
$ws.Range("K8").Value = 16000
$ws.Range("L8").Value = 17000
$ws.Range("M8").Value = 16500
$ws.Range("N8").Value = "$/caja 30 unidades"
$ws.Range("O8").Value = "Provincia de Limarí"
$ws.Range("P8").Value = 550
$ws.Range("Q8").Value = 30
$ws.Range("R8").Value = "Hortaliza"
